$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(205, 6).Value = 4938.666666666667
$ws.Cells.Item(205, 7).Value = 13656
$ws.Cells.Item(205, 8).Value = 325241.6666666667
$ws.Cells.Item(206, 6).Value = 5008.111111111111
$ws.Cells.Item(206, 7).Value = 13654.11111111111
$ws.Cells.Item(206, 8).Value = 323736.6666666667
$ws.Cells.Item(207, 6).Value = 4994.37037037037
$ws.Cells.Item(207, 7).Value = 13682.25925925926
$ws.Cells.Item(207, 8).Value = 322655.5555555556
$ws.Cells.Item(208, 6).Value = 4980.382716049383
$ws.Cells.Item(208, 7).Value = 13664.12345679012
$ws.Cells.Item(208, 8).Value = 323877.962962963
$ws.Cells.Item(209, 6).Value = 4994.288065843622
$ws.Cells.Item(209, 7).Value = 13666.83127572016
$ws.Cells.Item(209, 8).Value = 323423.3950617284
$ws.Cells.Item(210, 6).Value = 4989.680384087792
$ws.Cells.Item(210, 7).Value = 13671.07133058985
$ws.Cells.Item(210, 8).Value = 323318.9711934156
$ws.Cells.Item(211, 6).Value = 4988.117055326932
$ws.Cells.Item(211, 7).Value = 13667.34202103338
$ws.Cells.Item(211, 8).Value = 323540.109739369
$ws.Cells.Item(212, 6).Value = 4990.695168419448
$ws.Cells.Item(212, 7).Value = 13668.41487578113
$ws.Cells.Item(212, 8).Value = 323427.491998171
$ws.Cells.Item(213, 6).Value = 4989.497535944724
$ws.Cells.Item(213, 7).Value = 13668.94274246812
$ws.Cells.Item(213, 8).Value = 323428.8576436519
$ws.Cells.Item(214, 6).Value = 4989.436586563702
$ws.Cells.Item(214, 7).Value = 13668.23321309421
$ws.Cells.Item(214, 8).Value = 323465.4864603973
$ws.Cells.Item(215, 6).Value = 4989.876430309291
$ws.Cells.Item(215, 7).Value = 13668.53027711449
$ws.Cells.Item(215, 8).Value = 323440.6120340734
$ws.Cells.Item(216, 6).Value = 4989.603517605906
$ws.Cells.Item(216, 7).Value = 13668.56874422561
$ws.Cells.Item(216, 8).Value = 323444.9853793742
$ws.Cells.Item(217, 6).Value = 4989.6388448263
$ws.Cells.Item(217, 7).Value = 13668.44407814477
$ws.Cells.Item(217, 8).Value = 323450.3612912817
$ws.Cells.Item(218, 6).Value = 4989.706264247166
$ws.Cells.Item(218, 7).Value = 13668.51436649495
$ws.Cells.Item(218, 8).Value = 323445.3195682431
$ws.Cells.Item(219, 6).Value = 4989.649542226457
$ws.Cells.Item(219, 7).Value = 13668.50906295511
$ws.Cells.Item(219, 8).Value = 323446.8887462996
$ws.Cells.Item(220, 6).Value = 4989.664883766641
$ws.Cells.Item(220, 7).Value = 13668.48916919828
$ws.Cells.Item(220, 8).Value = 323447.5232019415
$ws.Cells.Item(221, 6).Value = 4989.673563413421
$ws.Cells.Item(221, 7).Value = 13668.50419954945
$ws.Cells.Item(221, 8).Value = 323446.5771721614
$ws.Cells.Item(222, 6).Value = 4989.662663135507
$ws.Cells.Item(222, 7).Value = 13668.50081056761
$ws.Cells.Item(222, 8).Value = 323446.9963734675
$ws.Cells.Item(223, 6).Value = 4989.667036771856
$ws.Cells.Item(223, 7).Value = 13668.49805977178
$ws.Cells.Item(223, 8).Value = 323447.0322491901
$ws.Cells.Item(224, 6).Value = 4989.667754440261
$ws.Cells.Item(224, 7).Value = 13668.50102329628
$ws.Cells.Item(224, 8).Value = 323446.868598273
$ws.Cells.Item(225, 6).Value = 4989.665818115875
$ws.Cells.Item(225, 7).Value = 13668.49996454522
$ws.Cells.Item(225, 8).Value = 323446.9657403102
$ws.Cells.Item(226, 6).Value = 4989.666869775997
$ws.Cells.Item(226, 7).Value = 13668.49968253776
$ws.Cells.Item(226, 8).Value = 323446.9555292578
$ws.Cells.Item(227, 6).Value = 4989.666814110711
$ws.Cells.Item(227, 7).Value = 13668.50022345975
$ws.Cells.Item(227, 8).Value = 323446.929955947
$ws.Cells.Item(228, 6).Value = 4989.666500667528
$ws.Cells.Item(228, 7).Value = 13668.49995684758
$ws.Cells.Item(228, 8).Value = 323446.950408505
$ws.Cells.Item(229, 6).Value = 4989.666728184745
$ws.Cells.Item(229, 7).Value = 13668.4999542817
$ws.Cells.Item(229, 8).Value = 323446.9452979033
$ws.Cells.Item(230, 6).Value = 4989.666680987662
$ws.Cells.Item(230, 7).Value = 13668.50004486301
$ws.Cells.Item(230, 8).Value = 323446.9418874517
$ws.Cells.Item(231, 6).Value = 4989.666636613311
$ws.Cells.Item(231, 7).Value = 13668.49998533076
$ws.Cells.Item(231, 8).Value = 323446.94586462
$ws.Cells.Item(232, 6).Value = 4989.666681928573
$ws.Cells.Item(232, 7).Value = 13668.49999482516
$ws.Cells.Item(232, 8).Value = 323446.9443499917
$ws.Cells.Item(233, 6).Value = 4989.666666509848
$ws.Cells.Item(233, 7).Value = 13668.50000833964
$ws.Cells.Item(233, 8).Value = 323446.9440340211
$ws.Cells.Item(234, 6).Value = 4989.666661683911
$ws.Cells.Item(234, 7).Value = 13668.49999616519
$ws.Cells.Item(234, 8).Value = 323446.9447495443
$ws.Cells.Item(235, 6).Value = 4989.666670040778
$ws.Cells.Item(235, 7).Value = 13668.49999977666
$ws.Cells.Item(235, 8).Value = 323446.9443778524
$ws.Cells.Item(236, 6).Value = 4989.666666078179
$ws.Cells.Item(236, 7).Value = 13668.50000142717
$ws.Cells.Item(236, 8).Value = 323446.9443871392
$ws.Cells.Item(237, 6).Value = 4989.666665934289
$ws.Cells.Item(237, 7).Value = 13668.49999912301
$ws.Cells.Item(237, 8).Value = 323446.9445048453
$ws.Cells.Item(238, 6).Value = 4989.666667351082
$ws.Cells.Item(238, 7).Value = 13668.50000010895
$ws.Cells.Item(238, 8).Value = 323446.944423279
$ws.Cells.Item(239, 6).Value = 4989.666666454516
$ws.Cells.Item(239, 7).Value = 13668.50000021971
$ws.Cells.Item(239, 8).Value = 323446.9444384212
$ws.Cells.Item(240, 6).Value = 4989.666666579962
$ws.Cells.Item(240, 7).Value = 13668.49999981722
$ws.Cells.Item(240, 8).Value = 323446.9444555151
$ws.Cells.Item(241, 6).Value = 4989.666666795187
$ws.Cells.Item(241, 7).Value = 13668.50000004862
$ws.Cells.Item(241, 8).Value = 323446.9444390718
$ws.Cells.Item(242, 6).Value = 4989.666666609888
$ws.Cells.Item(242, 7).Value = 13668.50000002852
$ws.Cells.Item(242, 8).Value = 323446.944444336
$ws.Cells.Item(243, 6).Value = 4989.666666661679
$ws.Cells.Item(243, 7).Value = 13668.49999996479
$ws.Cells.Item(243, 8).Value = 323446.9444463076
$ws.Cells.Item(244, 6).Value = 4989.666666688918
$ws.Cells.Item(244, 7).Value = 13668.50000001398
$ws.Cells.Item(244, 8).Value = 323446.9444432384
$ws.Cells.Item(245, 6).Value = 4989.666666653495
$ws.Cells.Item(245, 7).Value = 13668.50000000243
$ws.Cells.Item(245, 8).Value = 323446.9444446273
$ws.Cells.Item(246, 6).Value = 4989.66666666803
$ws.Cells.Item(246, 7).Value = 13668.49999999373
$ws.Cells.Item(246, 8).Value = 323446.9444447245
$ws.Cells.Item(247, 6).Value = 4989.666666670148
$ws.Cells.Item(247, 7).Value = 13668.50000000338
$ws.Cells.Item(247, 8).Value = 323446.9444441968
$ws.Cells.Item(248, 6).Value = 4989.666666663891
$ws.Cells.Item(248, 7).Value = 13668.49999999985
$ws.Cells.Item(248, 8).Value = 323446.9444445162
$ws.Cells.Item(249, 6).Value = 4989.666666667356
$ws.Cells.Item(249, 7).Value = 13668.49999999899
$ws.Cells.Item(249, 8).Value = 323446.9444444791
